$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-09-27"

# Update the column header text (B1 shared string) for the current-month column
$ws.Range("B1").Value = "September 2022 (through September 27)"

# Update/add the carjacking counts for newly reported data (added data for 2022-10-05)
$ws.Range("K2").Value = 11
$ws.Range("T2").Value = 11

$ws.Range("B3").Value = 8
$ws.Range("K3").Value = 16
$ws.Range("T3").Value = 7
$ws.Range("BM3").Value = 3

$ws.Range("B4").Value = 3

$ws.Range("BM6").Value = 4

$ws.Range("B7").Value = 2
$ws.Range("K7").Value = 6

$ws.Range("BD8").Value = 3

$ws.Range("BD9").Value = 5

$ws.Range("AU11").Value = 3

$ws.Range("K12").Value = 5

$ws.Range("K14").Value = 6

$ws.Range("K17").Value = 2

$ws.Range("T22").Value = 2

$ws.Range("K23").Value = 1
$ws.Range("AC23").Value = 4
$ws.Range("AL23").Value = 1

$ws.Range("B24").Value = 5
$ws.Range("K24").Value = 3

$ws.Range("AU25").Value = 1

$ws.Range("AU31").Value = 1

$ws.Range("BD37").Value = 1

$ws.Range("K38").Value = 6

$ws.Range("BD44").Value = 1

$ws.Range("K50").Value = 3

$ws.Range("K51").Value = 2

$ws.Range("K55").Value = 2
$ws.Range("AC55").Value = 1

$ws.Range("T57").Value = 4

$ws.Range("B58").Value = 1

$ws.Range("B64").Value = 4

$ws.Range("AU73").Value = 1

$ws.Range("AU78").Value = 1

$ws.Range("AC96").Value = 1
$ws.Range("BD96").Value = 2

$ws.Range("AC97").Value = 1

$ws.Range("T98").Value = 1
